# SAV-1028: Added support for ImagingArea product category
#
# Inserts two new reference-data sheets ("Imaging Type" and "X Ray Imaging
# Area") right after "Invoice Product", populates them, adds matching rows
# to the "Invoice Product" sheet, and renames the "Procedure" category
# value to "ProcedureType".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the two new sheets, in order, right after "Invoice Product".
# ---------------------------------------------------------------------------
$invoiceProduct = $wb.Worksheets.Item("Invoice Product")

$imagingType = $wb.Worksheets.Add($null, $invoiceProduct)
$imagingType.Name = "Imaging Type"

$xRayImagingArea = $wb.Worksheets.Add($null, $imagingType)
$xRayImagingArea.Name = "X Ray Imaging Area"

# ---------------------------------------------------------------------------
# 2. Populate "Imaging Type" with its header + single data row.
# ---------------------------------------------------------------------------
$imagingType.Cells.Item(1, 1).Value = "id"
$imagingType.Cells.Item(1, 2).Value = "code"
$imagingType.Cells.Item(1, 3).Value = "name"

$imagingType.Cells.Item(2, 1).Value = "ImagingType-xRay"
$imagingType.Cells.Item(2, 2).Value = "xRay"
$imagingType.Cells.Item(2, 3).Value = "X-Ray"

# ---------------------------------------------------------------------------
# 3. Populate "X Ray Imaging Area" with its header + single data row.
# ---------------------------------------------------------------------------
$xRayImagingArea.Cells.Item(1, 1).Value = "id"
$xRayImagingArea.Cells.Item(1, 2).Value = "code"
$xRayImagingArea.Cells.Item(1, 3).Value = "name"

$xRayImagingArea.Cells.Item(2, 1).Value = "xRayImagingArea-foot"
$xRayImagingArea.Cells.Item(2, 2).Value = "Xray-foot"
$xRayImagingArea.Cells.Item(2, 3).Value = "Foot"

# ---------------------------------------------------------------------------
# 4. Add the two new rows to "Invoice Product" describing the new products.
# ---------------------------------------------------------------------------
$invoiceProduct.Cells.Item(6, 1).Value = "InvoiceProduct-imagingType-xRay"
$invoiceProduct.Cells.Item(6, 2).Value = "ImagingType-xRay"
$invoiceProduct.Cells.Item(6, 3).Value = "X-Ray"
$invoiceProduct.Cells.Item(6, 4).Value = "ImagingType-xRay"
$invoiceProduct.Cells.Item(6, 5).Value = "ImagingType"
$invoiceProduct.Cells.Item(6, 6).Value = $true
$invoiceProduct.Cells.Item(6, 6).NumberFormat = '"TRUE";"TRUE";"FALSE"'
$invoiceProduct.Cells.Item(6, 7).Value = "current"

$invoiceProduct.Cells.Item(7, 1).Value = "InvoiceProduct-imagingArea-xRay-foot"
$invoiceProduct.Cells.Item(7, 2).Value = "ImagingArea-xRay-foot"
$invoiceProduct.Cells.Item(7, 3).Value = "X-Ray Foot"
$invoiceProduct.Cells.Item(7, 4).Value = "xRayImagingArea-foot"
$invoiceProduct.Cells.Item(7, 5).Value = "ImagingArea"
$invoiceProduct.Cells.Item(7, 6).Value = $true
$invoiceProduct.Cells.Item(7, 6).NumberFormat = '"TRUE";"TRUE";"FALSE"'
$invoiceProduct.Cells.Item(7, 7).Value = "current"

# ---------------------------------------------------------------------------
# 5. Rename the "Procedure" category value to "ProcedureType" (shared
#    string used by the Invoice Product row for procedure-34831).
# ---------------------------------------------------------------------------
$invoiceProduct.Cells.Item(4, 5).Value = "ProcedureType"

# ---------------------------------------------------------------------------
# 6. Restore tab/selection state: "X Ray Imaging Area" becomes the active
#    sheet (third tab, 0-indexed activeTab=2).
# ---------------------------------------------------------------------------
$invoiceProduct.Select()
$invoiceProduct.Range("D11").Select() | Out-Null
$xRayImagingArea.Select()
$xRayImagingArea.Range("A2").Select() | Out-Null
